# Narayan Jagadeesan .xlsx - add ownTeam/oppTeam columns and a new match row
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Non-breaking space, matching the trailing character already used after
# "Narayan Jagadeesan" elsewhere in this sheet (sheet name, F2/D2 cell).
$nbsp = [char]0x00A0

# Insert two new columns before the current "batsman" column (D),
# shifting batsman..sr from D:I to F:K.
$ws.Range("D:E").Insert()

# New header cells for the inserted columns.
$ws.Range("D1").Value = "ownTeam"
$ws.Range("E1").Value = "oppTeam"

# Fill in the new team columns for the existing data row (row 2).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "Chennai Super Kings"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "Mumbai Indians"

# Append a brand-new match row (row 3).
$ws.Range("A3").Value = " Dubai (DSC)"
$ws.Range("B3").Value = " October 10 2020"
$ws.Range("C3").Value = "RCB won by 37 runs"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "Chennai Super Kings"

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "Royal Challengers Bangalore"

$ws.Range("F3").Value = "Narayan Jagadeesan" + $nbsp

$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "33"

$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value = "28"

$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = "4"

$ws.Range("J3").NumberFormat = "@"
$ws.Range("J3").Value = "0"

$ws.Range("K3").NumberFormat = "@"
$ws.Range("K3").Value = "117.85"
